$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-F_Topic 0", "Politisk debat om indfødsret og statsborgerskab i Danmark"),
    @("2025-F_Topic 1", "Debat om statsborgerskab og integrationspolitik i Danmark"),
    @("2025-F_Topic 2", "Debatten om statsborgerskab og indvandringspolitik i Danmark"),
    @("2025-F_Topic 3", "Politisk Debat om Udlændingepolitik og Identitet i Danmark")
)

$startRow = 127
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
